# ---------------------------------------------------------------------------
# tc_p025v.docx edit — "create element rendition spec"
#
# Four content changes, all inside the transcription runs of this critical
# edition:
#   1. The <ms>pied</ms> element (third <ms>/</ms> pair in the document) is
#      wrapped in a new <bp> ... </bp> element.
#   2. One stray run of "</del>" was left with the wrong (blue, "<ms>"-style)
#      run color; it should carry the same brownish-red "<del>" color as
#      every other <del>/</del> run.
#   3. The word "costeaulx" (in "...seze costeaulx bien") gets wrapped with a
#      new <tl> ... </tl> element.
#   4. The word "foret" (in "...la bouche le foret passe uniment &") gets
#      wrapped with a new <tl> ... </tl> element.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Find the Nth occurrence of $searchText in the document and return it as a
# fresh Range (so further edits on it don't disturb the live Find range).
function Find-NthRange($searchText, $n) {
    $rng = $d.Content
    $count = 0
    while ($rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $count = $count + 1
        if ($count -eq $n) {
            return $d.Range($rng.Start, $rng.End)
        }
        $rng.Collapse(0)
        $rng.End = $d.Content.End
    }
    return $null
}

# Find the single run matching $searchText whose font color is $colorVal.
function Find-RangeByTextAndColor($searchText, $colorVal) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Font.Color = $colorVal
    if ($rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $true, "", 0)) {
        return $d.Range($rng.Start, $rng.End)
    }
    return $null
}

# --- Reference formatting for the blue Courier-New "tag" runs (<tl>, <ms>, ...) ---
# Captured up front, before any edits, from existing <tl> / </tl> runs so the
# brand-new tag runs we create end up with byte-for-byte the same run
# properties (rFonts ascii/eastAsia/hAnsi/cs, color, sz, szCs, rtl).
$refTagOpen = Find-NthRange "<tl>" 1
$refTagOpenFT = $refTagOpen.FormattedText

$refTagClose = Find-NthRange "</tl>" 1
$refTagCloseFT = $refTagClose.FormattedText

# ======================================================================
# Change 1 — wrap <ms>pied</ms> (the 3rd <ms>...</ms> pair) with <bp></bp>
# ======================================================================
$msOpen = Find-NthRange "<ms>" 3
$msOpen.Text = "<bp><ms>"

$msClose = Find-NthRange "</ms>" 3
$msClose.Text = "</ms></bp>"

# ======================================================================
# Change 2 — recolor the mis-colored "</del>" run (currently 0000ff / blue)
# to the standard deletion color a91111
# ======================================================================
$delClose = Find-RangeByTextAndColor "</del>" 16711680
$delClose.Font.Color = 1118633

# ======================================================================
# Change 3 — wrap "costeaulx" with <tl></tl>
# ======================================================================
$target3 = Find-NthRange "costeaulx" 1
$t3Start = $target3.Start
$t3Len = $target3.End - $target3.Start
$target3.Text = "<tl>costeaulx</tl>"

$tagOpen3 = $d.Range($t3Start, $t3Start + 4)
$tagOpen3.FormattedText = $refTagOpenFT

$tagClose3 = $d.Range($t3Start + 4 + $t3Len, $t3Start + 4 + $t3Len + 5)
$tagClose3.FormattedText = $refTagCloseFT

# ======================================================================
# Change 4 — wrap "foret" with <tl></tl>
# ======================================================================
$target4 = Find-NthRange "foret" 1
$t4Start = $target4.Start
$t4Len = $target4.End - $target4.Start
$target4.Text = "<tl>foret</tl>"

$tagOpen4 = $d.Range($t4Start, $t4Start + 4)
$tagOpen4.FormattedText = $refTagOpenFT

$tagClose4 = $d.Range($t4Start + 4 + $t4Len, $t4Start + 4 + $t4Len + 5)
$tagClose4.FormattedText = $refTagCloseFT

Write-Host "All changes applied."
